$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = 0.95833333333333337

$ws.Range("B13").Value = 43904
$ws.Range("C13").Value = 0.625
$ws.Range("D13").Value = 0.66666666666666663
$ws.Range("F13").Value = "CLion"
$ws.Range("G13").Value = "Bataille Navale"
$ws.Range("H13").Value = "Programmation du jeu"
$ws.Range("I13").Value = "Création des différentes fonctions"

$ws.Range("I15:I16").Select()
